# Se crea el servicio de delete, read y se cuadra los script correspondientes
# Update the testDataBooking sheet: keep 3 rows (simulate create/read/delete cycle),
# refresh booking data and remove the now-stale trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that are no longer part of the dataset (rows 5-11),
# bottom-up so the row indices of the rows above stay stable.
for ($r = 11; $r -ge 5; $r--) {
    $ws.Rows.Item($r).Delete()
}

# Row 2: booking was updated (and its deposit is no longer paid)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = "Cristhian"
$ws.Cells.Item(2, 4).Value = "Vargas"
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(2, 6).Value = '{"checkin":"2025-01-19","checkout":"2025-01-20"}'

# Row 3: new booking record
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = "Cristhian"
$ws.Cells.Item(3, 4).Value = "Montaño"
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = '{"checkin":"2025-01-19","checkout":"2025-01-20"}'

# Row 4: another new booking record
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = "David"
$ws.Cells.Item(4, 4).Value = "Vargas"
$ws.Cells.Item(4, 5).Value = $false
$ws.Cells.Item(4, 6).Value = '{"checkin":"2025-01-19","checkout":"2025-01-20"}'
